$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the underlying data values
$ws.Range("B6").Value = 2.67
$ws.Range("B7").Value = 113
$ws.Range("C14").Value = 2.33
$ws.Range("C15").Value = 1.67

# Update the SGPA formula (J17) to clamp values that round to ~2.0
$ws.Range("J17").Formula = "=IF(AND(ROUND(SUMPRODUCT(B12:B28,C12:C28)/SUM(B12:B28),2)>=1.996,ROUND(SUMPRODUCT(B12:B28,C12:C28)/SUM(B12:B28),2)<=1.999),2,ROUND(SUMPRODUCT(B12:B28,C12:C28)/SUM(B12:B28),2))"

# Update the CGPA formula (J18) the same way
$ws.Range("J18").Formula = "=IF(AND(ROUND(((SUMPRODUCT(B12:B28,C12:C28)) + ((B7*B6)-(SUMPRODUCT(F12:F26,G12:G26)))) / (((B7) - SUM(F12:F26)) + (SUM(B12:B26))),2) >= 1.996, ROUND(((SUMPRODUCT(B12:B28,C12:C28)) + ((B7*B6)-(SUMPRODUCT(F12:F26,G12:G26)))) / (((B7) - SUM(F12:F26)) + (SUM(B12:B26))),2) <= 1.999), 2, ROUND(((SUMPRODUCT(B12:B28,C12:C28)) + ((B7*B6)-(SUMPRODUCT(F12:F26,G12:G26)))) / (((B7) - SUM(F12:F26)) + (SUM(B12:B26))),2))"

# Clear the leftover formatted-but-empty cells that no longer carry any style
$ws.Range("J19:K19").Clear()
$ws.Range("I21:M21").Clear()
$ws.Range("J22:K22").Clear()

# Move the active selection to L23 as in the final saved state
$ws.Range("L23").Select()
